# Update cryptos list with latest price/volume data (auto-generated)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.449.88"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3
$ws.Range("D3").Value = "3.433.76"
$ws.Range("E3").Value = "  -4.34%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.06%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.70%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "135.28"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -8.72%  "

# Row 7
$ws.Range("D7").Value = "3.433.37"
$ws.Range("E7").Value = "  -4.33%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.15%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.490"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.61%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.43"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -5.68%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.121"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -10.63%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.377"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -9.01%  "

# Row 13
$ws.Range("D13").Value = "4.008.93"
$ws.Range("E13").Value = "  -4.54%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000180"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -12.49%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.47"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -10.65%  "

# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.442.76"
$ws.Range("E16").Value = "  -4.56%  "

# Row 17
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "65.271.05"
$ws.Range("E17").Value = "  -1.70%  "

# Row 18
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.114"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.92%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.01"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -9.85%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.73"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.44%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.66"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -8.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "392.56"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.14%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -6.70%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.544"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -11.08%  "

# Row 25
$ws.Range("E25").Value = "  -0.01%  "

# Row 26
$ws.Range("D26").Value = "3.573.33"
$ws.Range("E26").Value = "  -4.15%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000104"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -12.55%  "

# Row 28
$ws.Range("E28").Value = "  -0.17%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.25"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -9.65%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.17"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -13.15%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -12.81%  "

# Row 32
$ws.Range("D32").Value = "3.435.00"
$ws.Range("E32").Value = "  -4.24%  "

# Row 33
$ws.Range("E33").Value = "  +0.00%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.145"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.68%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "22.60"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -9.69%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "172.75"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.20%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -13.77%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.84"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -11.81%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.52"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -8.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.82"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -13.83%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0768"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -9.45%  "

# Row 42
$ws.Range("B42").Value = "Mantle"
$ws.Range("C42").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.811"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -8.33%  "

# Row 43
$ws.Range("B43").Value = "OKB"
$ws.Range("C43").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "43.66"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.82%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.41"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -15.04%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.62"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -12.32%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.09"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.55%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "21.95"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -6.40%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "6.54"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -8.42%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.13"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -15.49%  "

# Row 51
$ws.Range("D51").Value = "2.198.85"
$ws.Range("E51").Value = "  -7.92%  "

